$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.255.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.00%  "

# Row 3
$ws.Range("D3").Value = "'3.734.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'612.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.75%  "

# Row 6
$ws.Range("E6").Value = "  +11.19%  "

# Row 7
$ws.Range("D7").Value = "'0.641"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "

# Row 8
$ws.Range("E8").Value = "  -0.39%  "

# Row 9
$ws.Range("D9").Value = "'0.731"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.17%  "

# Row 10
$ws.Range("E10").Value = "  -1.30%  "

# Row 11
$ws.Range("D11").Value = "'60.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.35%  "

# Row 12
$ws.Range("E12").Value = "  -1.47%  "

# Row 13
$ws.Range("D13").Value = "'10.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.16%  "

# Row 14
$ws.Range("D14").Value = "'4.323.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.44%  "

# Row 15
$ws.Range("D15").Value = "'3.732.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
$ws.Range("D16").Value = "'1.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.72%  "

# Row 17
$ws.Range("D17").Value = "'19.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").Value = "'13.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "

# Row 20
$ws.Range("D20").Value = "'69.078.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "

# Row 21
$ws.Range("D21").Value = "'414.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.89%  "

# Row 22
$ws.Range("E22").Value = "  +2.21%  "

# Row 23
$ws.Range("D23").Value = "'90.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.39%  "

# Row 24
$ws.Range("D24").Value = "'3.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "

# Row 25
$ws.Range("D25").Value = "'11.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.83%  "

# Row 26
$ws.Range("D26").Value = "'13.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.02%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'3.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'6.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "

# Row 29
$ws.Range("D29").Value = "'9.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.50%  "

# Row 30
$ws.Range("D30").Value = "'33.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "

# Row 31
$ws.Range("E31").Value = "  +2.75%  "

# Row 32
$ws.Range("E32").Value = "  +1.80%  "

# Row 33
$ws.Range("D33").Value = "'651.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.66%  "

# Row 34
$ws.Range("D34").Value = "'0.124"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.69%  "

# Row 35
$ws.Range("D35").Value = "'46.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.40%  "

# Row 36
$ws.Range("D36").Value = "'66.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.77%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "'0.0₃0842"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.33%  "

# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.420"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.30%  "

# Row 39
$ws.Range("E39").Value = "  -0.10%  "

# Row 40
$ws.Range("E40").Value = "  +0.09%  "

# Row 41
$ws.Range("E41").Value = "  +4.64%  "

# Row 42
$ws.Range("D42").Value = "'3.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.38%  "

# Row 43
$ws.Range("E43").Value = "  +2.18%  "

# Row 44
$ws.Range("D44").Value = "'2.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.92%  "

# Row 45
$ws.Range("D45").Value = "'2.913.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.13%  "

# Row 46
$ws.Range("D46").Value = "'0.141"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'9.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "

# Row 48
$ws.Range("E48").Value = "  +0.99%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'144.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.34%  "

# Row 50
$ws.Range("E50").Value = "  -10.83%  "

# Row 51
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'3.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
